$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 128, shifting rows 128:140 down to 130:142
$ws.Range("A128:T129").EntireRow.Insert()

# New row 128 data
$ws.Cells.Item(128, 1).Value = 7
$ws.Cells.Item(128, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(128, 3).Value = "Ñuble"
$ws.Cells.Item(128, 4).Value = 44918
$ws.Cells.Item(128, 5).Value = 16
$ws.Cells.Item(128, 6).Value = "Fruta"
$ws.Cells.Item(128, 7).Value = 100103
$ws.Cells.Item(128, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(128, 9).Value = 100103001
$ws.Cells.Item(128, 10).Value = "Cereza"
$ws.Cells.Item(128, 11).Value = "Lapins"
$ws.Cells.Item(128, 12).Value = "Primera"
$ws.Cells.Item(128, 13).Value = 120
$ws.Cells.Item(128, 14).Value = 4000
$ws.Cells.Item(128, 15).Value = 5000
$ws.Cells.Item(128, 16).Value = 4500
$ws.Cells.Item(128, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(128, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(128, 19).Value = 450
$ws.Cells.Item(128, 20).Value = 10

# New row 129 data
$ws.Cells.Item(129, 1).Value = 7
$ws.Cells.Item(129, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(129, 3).Value = "Ñuble"
$ws.Cells.Item(129, 4).Value = 44918
$ws.Cells.Item(129, 5).Value = 16
$ws.Cells.Item(129, 6).Value = "Fruta"
$ws.Cells.Item(129, 7).Value = 100103
$ws.Cells.Item(129, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(129, 9).Value = 100103001
$ws.Cells.Item(129, 10).Value = "Cereza"
$ws.Cells.Item(129, 11).Value = "Lapins"
$ws.Cells.Item(129, 12).Value = "Segunda"
$ws.Cells.Item(129, 13).Value = 120
$ws.Cells.Item(129, 14).Value = 3000
$ws.Cells.Item(129, 15).Value = 3500
$ws.Cells.Item(129, 16).Value = 3250
$ws.Cells.Item(129, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(129, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(129, 19).Value = 325
$ws.Cells.Item(129, 20).Value = 10
